$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.038.99'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.40%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.264.94'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.94'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '183.90'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.45%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("E9").Value = '  -1.97%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.62'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.54%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.407'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.94%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.829.85'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.26%  '
$ws.Range("E13").Value = '  +1.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '68.027.27'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.55%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.33'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000168'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.262.73'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.71'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.47%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.24'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.94%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '416.28'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.51'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.27%  '
$ws.Range("E22").Value = '  +0.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.01'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.70%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.507'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.40%  '
$ws.Range("E25").Value = '  -2.74%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.188'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.29'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.35%  '
$ws.Range("E29").Value = '  -1.92%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.60'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.64%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.41'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.83%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.83'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.24'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.16%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '164.74'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.96%  '
$ws.Range("E35").Value = '  -5.18%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.88'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.68'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.791'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.93%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.43'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.59%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.27'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.626.59'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0673'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.48%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.39'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.46%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '333.89'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.44%  '
$ws.Range("E45").Value = '  -5.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0273'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.49%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.22'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.982'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0999'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.89%  '
$ws.Range("E50").Value = '  -0.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '30.54'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.47%  '
